# Update cryptos list with refreshed price/volume data from the latest
# coinranking.com scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column cells whose new values look like plain numbers need an
# explicit Text format first so Excel keeps them as strings (matching the
# source data, e.g. "372.93" must stay text, not become 372.93 as a number).
$textPriceCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D12",
    "D14",
    "D15",
    "D16",
    "D18",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D32",
    "D33",
    "D35",
    "D37",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D49",
    "D51",
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "50.665.23"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.915.86"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("D5").Value = "372.93"
$ws.Range("E5").Value = "  -3.04%  "
$ws.Range("D6").Value = "99.44"
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").Value = "35.46"
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "0.0842"
$ws.Range("D13").Value = "3.371.07"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").Value = "17.88"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").Value = "11.97"
$ws.Range("E15").Value = "  +62.11%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "7.50"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "2.923.72"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "0.983"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "50.646.22"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("E20").Value = "  -8.43%  "
$ws.Range("D21").Value = "12.14"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").Value = "0.0₃0939"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "68.97"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "265.03"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "3.10"
$ws.Range("E25").Value = "  +6.91%  "
$ws.Range("D26").Value = "7.77"
$ws.Range("E26").Value = "  -5.06%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "7.08"
$ws.Range("E28").Value = "  -6.41%  "
$ws.Range("D29").Value = "25.29"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  -5.27%  "
$ws.Range("E31").Value = "  -5.69%  "
$ws.Range("D32").Value = "9.83"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "50.23"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "32.86"
$ws.Range("E35").Value = "  -4.65%  "
$ws.Range("D37").Value = "0.0426"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "16.20"
$ws.Range("E40").Value = "  -5.02%  "
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").Value = "2.40"
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("D43").Value = "119.05"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("D44").Value = "20.76"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "2.02"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "1.977.13"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").Value = "0.256"
$ws.Range("E49").Value = "  -7.71%  "
$ws.Range("E50").Value = "  -7.00%  "
$ws.Range("D51").Value = "5.21"
